$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Province" header (A1) becomes "* Province" to indicate a required field.
$ws.Range("A1").Value = "* Province"
